# "Drop in files from RMI script"
# This script reproduces dropping in an updated version of the workbook that:
#  - no longer carries the Texas-specific analysis ("Texas Notes" sheet)
#  - re-points the EoDSDwSP outputs straight at the national "Calculations" sheet
#  - removes the now-unused hyperlink on the About sheet

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook

$wsAbout        = $wb.Worksheets.Item("About")
$wsCalculations = $wb.Worksheets.Item("Calculations")
$wsTexasNotes   = $wb.Worksheets.Item("Texas Notes")
$wsEoDSDwSP     = $wb.Worksheets.Item("EoDSDwSP")

# --- EoDSDwSP: stop relying on the Texas Notes tab, reference Calculations directly ---
$wsEoDSDwSP.Range("B2").Formula = "=Calculations!B9"
$wsEoDSDwSP.Range("B4").Formula = "=Calculations!B10"

# update the on-sheet selection before we move off of it
$wsEoDSDwSP.Range("B2").Select()

# --- About: drop the hyperlink that pointed into the source PDF reference cell ---
$wsAbout.Range("B6").Style = "Normal"
$wsAbout.Hyperlinks.Delete()

# the Hyperlink cell style is no longer used anywhere in the workbook
$wb.Styles.Item("Hyperlink").Delete()

# --- Calculations: reset selection back to the top-left of the sheet ---
$wsCalculations.Activate()
$wsCalculations.Range("A1").Select()

# --- Remove the Texas Notes sheet entirely ---
$wsTexasNotes.Delete()

# --- About becomes the active sheet/selection when the workbook is opened ---
$wsAbout.Activate()
$wsAbout.Range("A12").Select()
